$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# (Target widths from the source file use fractional-pixel metrics that this
# engine's ColumnWidth setter quantizes to 1/7-character steps; the values
# below are chosen so the stored width lands as close as possible to the
# target after that quantization.)
$ws.Columns.Item(14).ColumnWidth = 4.571428571428571   # N -> stored 5.285714285714286 (target 5.289887640449439)
$ws.Columns.Item(15).ColumnWidth = 5.714285714285714   # O -> stored 6.428571428571429 (target 6.389887640449439)
$ws.Columns.Item(21).ColumnWidth = 7.857142857142857   # U -> stored 8.571428571428571 (target 8.589887640449438)
$ws.Columns.Item(22).ColumnWidth = 7.857142857142857   # V -> stored 8.571428571428571 (target 8.589887640449438)

# --- Row 2 (Sekhar Beri) ---
$ws.Range("J2").Value = 10000.0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 6000.0
$ws.Range("R2").Value = 10000.0
$ws.Range("S2").Value = 0.0
$ws.Range("T2").Value = 0.0
$ws.Range("Y2").Value = 0.0
$ws.Range("Z2").Value = 10000.0

# --- Row 3 (Priyanka Muddana) ---
$ws.Range("J3").Value = 9345.0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 5345.0
$ws.Range("R3").Value = 9345.0
$ws.Range("T3").Value = 163.54
$ws.Range("Y3").Value = 643.54
$ws.Range("Z3").Value = 8701.46

# --- Row 4 (Pattabhi RamaRao Galidevara) ---
$ws.Range("J4").Value = 9520.0
$ws.Range("P4").Value = 5520.0
$ws.Range("R4").Value = 9520.0
$ws.Range("S4").Value = 480.0
$ws.Range("U4").Value = 1200.0
$ws.Range("V4").Value = 1400.0
$ws.Range("Y4").Value = 3080.0
$ws.Range("Z4").Value = 6440.0
